# Append a new paragraph at the end of the document body, after the
# paragraph that ends with "difference in velocity."
#
# The new paragraph contains four runs (all sharing identical run
# formatting) whose concatenated text reads:
#   "Wijnand’s model assumes that each particle can only undergo one
#   collision each timestep, no matter how large the time step. This
#   assumption is only holds if the spatial and temporal resolutions are
#   small enough that the collision rate is always smaller than one."
# and the first run carries a <w:lastRenderedPageBreak/> marker before its
# text, matching the paragraph break Word recorded when it last laid the
# document out across pages.

$d = $word.ActiveDocument

# Move to the very end of the document and open a brand-new paragraph
# there; Word seeds its paragraph mark (pPr/rPr) from the paragraph it
# follows, so it automatically inherits the "jc=both" / eastAsia font /
# en-GB language formatting already used throughout this section.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# Grab the freshly created (still empty) paragraph and drop the four runs
# of OOXML straight into it via Range.InsertXML so that the
# <w:lastRenderedPageBreak/> marker -- which has no dedicated COM
# property -- lands in the exact spot Word itself would put it.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

$fragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Wijnand’s model assumes that each particle can only undergo one collision each timestep, no matter how large the time step. This assumption is </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">only holds </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:val="en-GB"/></w:rPr><w:t>if the spatial and temporal resolutions are small enough that the collision rate is always smaller than one</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:val="en-GB"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@

$newRange.InsertXML($fragment) | Out-Null

Write-Output "Inserted paragraph; document now has $($d.Paragraphs.Count) paragraphs."
